$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New rows 9-11: names first (column A) ---
$ws.Range("A9").Value = "RamKumar"
$ws.Range("A10").Value = "Sharath"
$ws.Range("A11").Value = "Senthil"

$ws.Range("B9").Value = 9
$ws.Range("B10").Value = 10
$ws.Range("B11").Value = 12

# --- New category column C for rows 1-11 ---
$ws.Range("C1").Value = "ADK01CBE16"
$ws.Range("C2").Value = "ADK02CBE16"
$ws.Range("C3").Value = "ADK03CBE16"
$ws.Range("C4").Value = "ADK04CBE16"
$ws.Range("C6").Value = "ADK06CBE16"
$ws.Range("C5").Value = "ADK05CBE16"
$ws.Range("C8").Value = "ADK08CBE16"
$ws.Range("C7").Value = "ADK07CBE16"
$ws.Range("C9").Value = "ADK09CBE16"
$ws.Range("C10").Value = "ADK10CBE16"
$ws.Range("C11").Value = "ADK11CBE16"

# --- New rows 12-16 (name, category, seeding) ---
$ws.Range("A12").Value = "Subbu"
$ws.Range("C12").Value = "ADK12CBE16"
$ws.Range("B12").Value = 13

$ws.Range("A13").Value = "Sennu"
$ws.Range("C13").Value = "ADK13CBE16"
$ws.Range("B13").Value = 16

$ws.Range("A14").Value = "Ganesh"
$ws.Range("C14").Value = "ADK14CBE16"
$ws.Range("B14").Value = 18

$ws.Range("A15").Value = "Jagan"
$ws.Range("C15").Value = "ADK15CBE16"
$ws.Range("B15").Value = 20

$ws.Range("A16").Value = "Naren"
$ws.Range("C16").Value = "ADK16CBE16"
$ws.Range("B16").Value = 19

# Column C width (target stored width 12.28515625 chars; closest attainable
# via the COM ColumnWidth->pixel rounding is ColumnWidth=11.5 -> 12.3333...)
$ws.Columns.Item(3).ColumnWidth = 11.5

# Selection moves to A17
$ws.Range("A17").Select()
